# Daily attendance processing - 2026-01-07 06:04:59
# Normalize the "Recorded By" (column G) entries so that the actor who
# performed the recording is listed first and "System" is no longer the
# leading entry: swap the first two comma-separated names whenever the
# first one is the literal "System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value()

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "

        if ($parts.Length -gt 1 -and $parts[0] -eq "System") {
            $newParts = New-Object System.Collections.ArrayList
            [void]$newParts.Add($parts[1])
            [void]$newParts.Add($parts[0])
            for ($i = 2; $i -lt $parts.Length; $i++) {
                [void]$newParts.Add($parts[$i])
            }
            $cell.Value = ($newParts -join ", ")
        }
    }
}
